$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" list:
# "LOM3246: Técnicas de Caracterização de Materiais (Requisito)"
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "LOM3246: Técnicas de Caracterização de Materiais (Requisito)",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the LOM3246 requisito paragraph"
}

# Identify that paragraph's index within the document's Paragraphs collection.
$allParas = $d.Paragraphs
$lomIndex = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {
        $lomIndex = $i
    }
}

if ($lomIndex -eq -1) {
    throw "Could not resolve paragraph index for LOM3246 requisito paragraph"
}

# The three paragraphs right after it are:
#   1) an empty separator paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# All three (including their paragraph marks) are removed, leaving the
# trailing empty paragraph that used to follow them immediately after the
# LOM3246 requisito paragraph.
$firstToRemove = $allParas.Item($lomIndex + 1)
$lastToRemove = $allParas.Item($lomIndex + 3)

$deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$deleteRange.Delete()
